$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain-text numbers (e.g. "116.78"); some new
# values would otherwise be auto-converted to floating point numbers by
# Excel, so force those specific cells to Text format first.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range('D2').Value = '43.798.36'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '2.291.85'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '116.78'
$ws.Range('E5').Value = '  +13.24%  '
$ws.Range('D6').Value = '269.46'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').Value = '49.41'
$ws.Range('E10').Value = '  +7.70%  '
$ws.Range('D11').Value = '0.0947'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '9.01'
$ws.Range('E12').Value = '  +13.35%  '
$ws.Range('D13').Value = '0.108'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = '15.86'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '2.633.95'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').Value = '0.878'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').Value = '2.288.10'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').Value = '43.678.48'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '6.99'
$ws.Range('E20').Value = '  +11.80%  '
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').Value = '2.43'
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('D23').Value = '10.07'
$ws.Range('E23').Value = '  +9.40%  '
$ws.Range('D24').Value = '233.65'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '11.72'
$ws.Range('E27').Value = '  +4.16%  '
$ws.Range('D28').Value = '3.93'
$ws.Range('E28').Value = '  +2.64%  '
$ws.Range('D29').Value = '42.05'
$ws.Range('E29').Value = '  +8.61%  '
$ws.Range('E30').Value = '  -2.07%  '
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').Value = '173.86'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').Value = '0.0938'
$ws.Range('E33').Value = '  +4.68%  '
$ws.Range('D34').Value = '21.61'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('E35').Value = '  +4.47%  '
$ws.Range('D36').Value = '0.128'
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').Value = '4.74'
$ws.Range('E37').Value = '  -0.86%  '
$ws.Range('D38').Value = '0.0361'
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('E40').Value = '  +7.35%  '
$ws.Range('D41').Value = '14.65'
$ws.Range('E41').Value = '  +19.62%  '
$ws.Range('D42').Value = '74.75'
$ws.Range('E42').Value = '  +15.30%  '
$ws.Range('E43').Value = '  +3.81%  '
$ws.Range('D44').Value = '0.243'
$ws.Range('E44').Value = '  +3.27%  '
$ws.Range('D45').Value = '6.39'
$ws.Range('E45').Value = '  +21.62%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('D48').Value = '8.77'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('D49').Value = '103.01'
$ws.Range('E49').Value = '  +4.40%  '
$ws.Range('E50').Value = '  +3.78%  '
$ws.Range('E51').Value = '  -1.96%  '
